# ---------------------------------------------------------------------------
# Rename "Intervention coverages" -> "Interventions coverages" and add three
# new sheets describing the interventions (affected fraction, mortality
# effectiveness, incidence effectiveness), matching the target OOXML diff.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$covWs = $wb.Worksheets.Item("Intervention coverages")
$covWs.Name = "Interventions coverages"

# New intervention rows to append to the coverages sheet (rows 4-10).
$newInterventions = @(
    "Complementary feeding 1",
    "Complementary feeding 2",
    "Complementary feeding 3",
    "Breastfeeding promotion",
    "IPTp",
    "BES",
    "MMS"
)

$row = 4
foreach ($name in $newInterventions) {
    $covWs.Cells.Item($row, 1).Value = $name
    $covWs.Cells.Item($row, 2).Value = 0.0
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# Helper data used by the three new sheets: each row is
# (intervention label shown once per pair of rows, cause, values C..G)
# ---------------------------------------------------------------------------

$headers = @("Interventions", "Cause", "<1 month", "1-5 months", "6-11 months", "12-23 months", "24-59 months")

function Fill-InterventionSheet {
    param($ws, $f25923, $f251, $f416, $f300)

    $ws.Cells.Item(1, 1).Value = $headers[0]
    $ws.Cells.Item(1, 2).Value = $headers[1]
    $ws.Cells.Item(1, 3).Value = $headers[2]
    $ws.Cells.Item(1, 4).Value = $headers[3]
    $ws.Cells.Item(1, 5).Value = $headers[4]
    $ws.Cells.Item(1, 6).Value = $headers[5]
    $ws.Cells.Item(1, 7).Value = $headers[6]

    $ws.Cells.Item(2, 1).Value = "Zinc supplementation"
    $ws.Cells.Item(2, 2).Value = "Diarrhea"
    $ws.Cells.Item(2, 3).Value = 0.0
    $ws.Cells.Item(2, 4).Value = 0.0
    $ws.Cells.Item(2, 5).Value = 0.0
    $ws.Cells.Item(2, 6).Value = $f25923[0]
    $ws.Cells.Item(2, 7).Value = $f25923[0]

    $ws.Cells.Item(3, 2).Value = "Pneumonia"
    $ws.Cells.Item(3, 3).Value = 0.0
    $ws.Cells.Item(3, 4).Value = 0.0
    $ws.Cells.Item(3, 5).Value = 0.0
    $ws.Cells.Item(3, 6).Value = $f25923[1]
    $ws.Cells.Item(3, 7).Value = $f25923[1]

    $ws.Cells.Item(4, 1).Value = "Vitamin A supplementation"
    $ws.Cells.Item(4, 2).Value = "Diarrhea"
    $ws.Cells.Item(4, 3).Value = 0.0
    $ws.Cells.Item(4, 4).Value = 0.0
    $ws.Cells.Item(4, 5).Value = $f416[0]
    $ws.Cells.Item(4, 6).Value = $f416[0]
    $ws.Cells.Item(4, 7).Value = $f416[0]

    $ws.Cells.Item(5, 2).Value = "Pneumonia"
    $ws.Cells.Item(5, 3).Value = 0.0
    $ws.Cells.Item(5, 4).Value = 0.0
    $ws.Cells.Item(5, 5).Value = 0.0
    $ws.Cells.Item(5, 6).Value = 0.0
    $ws.Cells.Item(5, 7).Value = 0.0
}

# ---------------------------------------------------------------------------
# "Interventions affected fraction"
# ---------------------------------------------------------------------------
$affWs = $wb.Worksheets.Add($null, $covWs)
$affWs.Name = "Interventions affected fraction"
Fill-InterventionSheet $affWs @(0.253, 0.253) @(0.5, 0.5) @(0.416, 0.416) @(0.3, 0.3)

# ---------------------------------------------------------------------------
# "Interventions mortality eff"
# ---------------------------------------------------------------------------
$mortWs = $wb.Worksheets.Add($null, $affWs)
$mortWs.Name = "Interventions mortality eff"
$mortWs.Cells.Item(1, 1).Value = $headers[0]
$mortWs.Cells.Item(1, 2).Value = $headers[1]
$mortWs.Cells.Item(1, 3).Value = $headers[2]
$mortWs.Cells.Item(1, 4).Value = $headers[3]
$mortWs.Cells.Item(1, 5).Value = $headers[4]
$mortWs.Cells.Item(1, 6).Value = $headers[5]
$mortWs.Cells.Item(1, 7).Value = $headers[6]

$mortWs.Cells.Item(2, 1).Value = "Zinc supplementation"
$mortWs.Cells.Item(2, 2).Value = "Diarrhea"
$mortWs.Cells.Item(2, 3).Value = 0.0
$mortWs.Cells.Item(2, 4).Value = 0.0
$mortWs.Cells.Item(2, 5).Value = 0.0
$mortWs.Cells.Item(2, 6).Value = 0.5
$mortWs.Cells.Item(2, 7).Value = 0.5

$mortWs.Cells.Item(3, 2).Value = "Pneumonia"
$mortWs.Cells.Item(3, 3).Value = 0.0
$mortWs.Cells.Item(3, 4).Value = 0.0
$mortWs.Cells.Item(3, 5).Value = 0.0
$mortWs.Cells.Item(3, 6).Value = 0.51
$mortWs.Cells.Item(3, 7).Value = 0.51

$mortWs.Cells.Item(4, 1).Value = "Vitamin A supplementation"
$mortWs.Cells.Item(4, 2).Value = "Diarrhea"
$mortWs.Cells.Item(4, 3).Value = 0.0
$mortWs.Cells.Item(4, 4).Value = 0.0
$mortWs.Cells.Item(4, 5).Value = 0.3
$mortWs.Cells.Item(4, 6).Value = 0.3
$mortWs.Cells.Item(4, 7).Value = 0.3

$mortWs.Cells.Item(5, 2).Value = "Pneumonia"
$mortWs.Cells.Item(5, 3).Value = 0.0
$mortWs.Cells.Item(5, 4).Value = 0.0
$mortWs.Cells.Item(5, 5).Value = 0.0
$mortWs.Cells.Item(5, 6).Value = 0.0
$mortWs.Cells.Item(5, 7).Value = 0.0

# ---------------------------------------------------------------------------
# "Interventions incidence eff"
# ---------------------------------------------------------------------------
$incWs = $wb.Worksheets.Add($null, $mortWs)
$incWs.Name = "Interventions incidence eff"
$incWs.Cells.Item(1, 1).Value = $headers[0]
$incWs.Cells.Item(1, 2).Value = $headers[1]
$incWs.Cells.Item(1, 3).Value = $headers[2]
$incWs.Cells.Item(1, 4).Value = $headers[3]
$incWs.Cells.Item(1, 5).Value = $headers[4]
$incWs.Cells.Item(1, 6).Value = $headers[5]
$incWs.Cells.Item(1, 7).Value = $headers[6]

$incWs.Cells.Item(2, 1).Value = "Zinc supplementation"
$incWs.Cells.Item(2, 2).Value = "Diarrhea"
$incWs.Cells.Item(2, 3).Value = 0.0
$incWs.Cells.Item(2, 4).Value = 0.0
$incWs.Cells.Item(2, 5).Value = 0.0
$incWs.Cells.Item(2, 6).Value = 0.65
$incWs.Cells.Item(2, 7).Value = 0.65

$incWs.Cells.Item(3, 2).Value = "Pneumonia"
$incWs.Cells.Item(3, 3).Value = 0.0
$incWs.Cells.Item(3, 4).Value = 0.0
$incWs.Cells.Item(3, 5).Value = 0.0
$incWs.Cells.Item(3, 6).Value = 0.52
$incWs.Cells.Item(3, 7).Value = 0.52

$incWs.Cells.Item(4, 1).Value = "Vitamin A supplementation"
$incWs.Cells.Item(4, 2).Value = "Diarrhea"
$incWs.Cells.Item(4, 3).Value = 0.0
$incWs.Cells.Item(4, 4).Value = 0.0
$incWs.Cells.Item(4, 5).Value = 0.62
$incWs.Cells.Item(4, 6).Value = 0.62
$incWs.Cells.Item(4, 7).Value = 0.62

$incWs.Cells.Item(5, 2).Value = "Pneumonia"
$incWs.Cells.Item(5, 3).Value = 0.0
$incWs.Cells.Item(5, 4).Value = 0.0
$incWs.Cells.Item(5, 5).Value = 0.0
$incWs.Cells.Item(5, 6).Value = 0.0
$incWs.Cells.Item(5, 7).Value = 0.0

Write-Host "Done. Worksheets:" $wb.Worksheets.Count
